$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_cases")

# 1. Remove the existing autofilter criterion (the "login page" filter) and
#    unhide the rows it was hiding (rows 3-5). This clears <filterColumn>
#    and the hidden="1" flags on those rows.
$ws.ShowAllData()

# 2. Add the new "profile page" test-case row (row 11). Values are entered
#    in the same order the original author would have typed them (A, B,
#    then D before C) so the shared-string table gets the same new-entry
#    ordering as the target workbook.
$ws.Range("A11").Value = "profile page"
$ws.Range("B11").Value = "profile page test"
$ws.Range("D11").Value = "Editing the profile"
$ws.Range("C11").Value = "test_changing_profile"
$ws.Range("F11").Value = "NORMAL"
$ws.Range("G11").Value = "NO"

# 3. Grow the autofilter range so it covers the new row. Calling
#    Range.AutoFilter() while autofilter is already active just toggles it
#    off (classic VBA behaviour), so switch it off explicitly first and
#    then reapply it over the bigger range.
$ws.AutoFilterMode = $False
$ws.Range("A2:H11").AutoFilter()

# 4. Keep the workbook-level hidden _FilterDatabase name in sync with the
#    resized autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "test_cases!_FilterDatabase") {
        $n.RefersTo = "=test_cases!`$A`$2:`$H`$11"
    }
}

# 5. Match the author's final selection (cell C11, where they'd just typed
#    the new test-case name).
$ws.Range("C11").Select()
